$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "N摩尔-U"
$ws.Range("C2").Value = "摩尔线程"
$ws.Range("A3").Value = "N摩尔-U"
$ws.Range("B3").Value = "航天发展"
$ws.Range("C3").Value = "航天发展"
$ws.Range("A4").Value = "合富中国"
$ws.Range("B4").Value = "平潭发展"
$ws.Range("C4").Value = "平潭发展"
$ws.Range("A5").Value = "平潭发展"
$ws.Range("B5").Value = "和而泰"
$ws.Range("C5").Value = "海王生物"
$ws.Range("A6").Value = "航天机电"
$ws.Range("B6").Value = "合富中国"
$ws.Range("C6").Value = "合富中国"
$ws.Range("A7").Value = "航天动力"
$ws.Range("B7").Value = "海王生物"
$ws.Range("A8").Value = "龙洲股份"
$ws.Range("C8").Value = "航天动力"
$ws.Range("A9").Value = "海王生物"
$ws.Range("B9").Value = "实达集团"
$ws.Range("C9").Value = "和而泰"
$ws.Range("A10").Value = "和而泰"
$ws.Range("B10").Value = "巨轮智能"
$ws.Range("C10").Value = "龙洲股份"
$ws.Range("B11").Value = "龙洲股份"
$ws.Range("C11").Value = "顺灏股份"
$ws.Range("A12").Value = "乾照光电"
$ws.Range("B12").Value = "东方财富"
$ws.Range("C12").Value = "海欣食品"
$ws.Range("A13").Value = "实达集团"
$ws.Range("B13").Value = "航天动力"
$ws.Range("C13").Value = "航天机电"
$ws.Range("A14").Value = "顺灏股份"
$ws.Range("B14").Value = "永鼎股份"
$ws.Range("C14").Value = "太阳电缆"
$ws.Range("A15").Value = "永鼎股份"
$ws.Range("B15").Value = "海欣食品"
$ws.Range("C15").Value = "永鼎股份"
$ws.Range("A16").Value = "赢时胜"
$ws.Range("B16").Value = "中银证券"
$ws.Range("C16").Value = "特发信息"
$ws.Range("A17").Value = "雷科防务"
$ws.Range("B17").Value = "赢时胜"
$ws.Range("C17").Value = "赢时胜"
$ws.Range("A18").Value = "巨轮智能"
$ws.Range("B18").Value = "雷科防务"
$ws.Range("A19").Value = "超捷股份"
$ws.Range("B19").Value = "乾照光电"
$ws.Range("C19").Value = "乾照光电"
$ws.Range("A20").Value = "中银证券"
$ws.Range("B20").Value = "航天科技"
$ws.Range("C20").Value = "安记食品"
$ws.Range("A21").Value = "特发信息"
$ws.Range("B21").Value = "三花智控"
$ws.Range("C21").Value = "雷科防务"
